$d = $word.ActiveDocument

# Change 1: merge the "(incluindo você)" runs into the sentence as plain text.
$d.Content.Find.Execute(
    "O que você achou que mandamos bem na Jornada (incluindo você); ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "O que você achou que mandamos bem na Jornada (incluindo você); ", 0)

# Change 2: "dos bugs ... os ciclos" -> "do erro ... o ciclo"
$d.Content.Find.Execute(
    "Poderia melhorar a questão dos bugs ocorridos durante os ciclos ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Poderia melhorar a questão do erro ocorridos durante o ciclo ", 2)

# Change 3: "mas era ao na aplicação. " -> "mas era na aplicação."
$d.Content.Find.Execute(
    "mas era ao na aplicação. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "mas era na aplicação.", 2)
